# Insert a new "Hype 炒作" vocabulary-list entry as its own paragraph,
# immediately after the existing "Prophet 預言家" paragraph and before the
# trailing (empty, bookmarked) paragraph at the end of the document.
#
# The new paragraph has three runs, matching the formatting convention used
# by the other two-letter-prefix entries in this document (e.g. "Halo 光環"):
#   1. "H"     - rFonts hint=eastAsia, no East-Asian language override
#   2. "ype "  - plain run, no run properties at all
#   3. "炒作"   - rFonts hint=eastAsia + lang eastAsia=zh-HK

$d = $word.ActiveDocument

# Find the "Prophet 預言家" paragraph -- the new entry goes right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*預言家*") {
        $target = $p
    }
}

# Insert a fresh empty paragraph right after it.
$null = $target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Populate the new paragraph with the exact run/formatting structure via a
# WordOpenXML fragment, so each run keeps precisely the run properties we
# want (rather than inheriting/merging with neighboring run formatting).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>H</w:t></w:r><w:r><w:t xml:space="preserve">ype </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-HK"/></w:rPr><w:t>炒作</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newPara.Range.InsertXML($xml)
